$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells (I1, J1) - set values first
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the formatting (bold, border, centered alignment) from the existing
# header cell H1 onto the two new header cells so they match the rest of
# row 1.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# New data cells (I2:J3)
$ws.Range("I2").Value = 5
$ws.Range("J2").Value = 7
$ws.Range("I3").Value = 7
$ws.Range("J3").Value = 9
